$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("SeriesInfo")
$ws2.Range("B3").Value = "'2023-12-10"
$ws2.Range("B3").Style = "Normal"
